$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (theta_se) - updated standard errors from new bootstrapping
$ws.Range("B4").Value = "(2.57)"
$ws.Range("C4").Value = "(0.62)"
$ws.Range("D4").Value = "(0.57)"
$ws.Range("E4").Value = "(2.71)"
$ws.Range("F4").Value = "(0.27)"
$ws.Range("G4").Value = "(0.57)"
$ws.Range("H4").Value = "(0.47)"
$ws.Range("I4").Value = "(2.09)"
$ws.Range("J4").Value = "(0.7)"
$ws.Range("K4").Value = "(2.0)"
$ws.Range("L4").Value = "(0.09)"

# Row 6 (lambda_se) - updated standard errors from new bootstrapping
$ws.Range("B6").Value = "(1.96)"
$ws.Range("C6").Value = "(0.1)"
$ws.Range("D6").Value = "(0.56)"
$ws.Range("E6").Value = "(1.46)"
$ws.Range("F6").Value = "(0.55)"
$ws.Range("G6").Value = "(0.09)"
$ws.Range("H6").Value = "(0.36)"
$ws.Range("I6").Value = "(1.45)"
$ws.Range("J6").Value = "(0.8)"
$ws.Range("K6").Value = "(0.51)"
$ws.Range("L6").Value = "(1.56)"
